$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.837.64"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.730.73"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4835"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "1.729.63"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06870"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6039"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "26.823.10"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007137"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "1.952.71"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.403"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.447"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.066"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.797"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.950"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07913"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.660"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04561"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.597"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6164"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9254"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.451"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.988"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.664"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.20%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01493"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3831"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.790"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1156"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05363"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.919"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.240"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
